# ProductBacklog.xlsx update "uppfaert skv. abendingum kennara"
# (updated according to teacher's comments)
#
# - Re-balance the "Sprint" (col A) and "Priority" (col D) values for the
#   backlog items in rows 4-15.
# - Clear the leftover placeholder "..." cells in F4:G4.
# - Add a new backlog item "Hotel-based pickup" in C17.
# - Cosmetic view changes: zoom to 200%, move the active cell to D6, and
#   let the default column width re-flow (Excel recomputes this on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Search day tours ---
$ws.Range("A4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("F4:G4").Clear()

# --- Row 5: Filter tours ---
$ws.Range("A5").Value = 1
$ws.Range("D5").Value = 1

# --- Row 6: Read reviews for tours (unchanged values, kept for clarity) ---
$ws.Range("A6").Value = 2
$ws.Range("D6").Value = 3

# --- Row 7: Set price range ---
$ws.Range("A7").Value = 1
$ws.Range("D7").Value = 1

# --- Row 8: See tour itinerary ---
$ws.Range("A8").Value = 1
$ws.Range("D8").Value = 2

# --- Row 9: Cancel tour ---
$ws.Range("A9").Value = 2
$ws.Range("D9").Value = 4

# --- Row 10: Search by attractions ---
$ws.Range("A10").Value = 1

# --- Row 11: Pay for tours ---
$ws.Range("A11").Value = 2
$ws.Range("D11").Value = 4

# --- Row 12: Email tour ticket ---
$ws.Range("A12").Value = 2

# --- Row 13: add/remove listings ---
$ws.Range("A13").Value = 2

# --- Row 14: Receive bookings ---
$ws.Range("A14").Value = 2

# --- Row 15: Set availability ---
$ws.Range("A15").Value = 2

# --- New backlog item ---
$ws.Range("C17").Value = "Hotel-based pickup"

# --- View / selection cosmetics ---
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.Zoom = 200
$null = $ws.Range("D6").Select()
